# Add fake data for database
# - rename "customer_status" column header to "cup_rental" and repurpose it to
#   hold a (mostly empty) numeric rental id
# - remove the stray "Borrow" value from the old customer_status column
# - bump row2's join_date
# - give every customer a "deposit" of 5 and populate "cups_bought" for each row
# - clear out the now-unused account_value data cell
# - move the active selection to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: customer_status -> cup_rental
$ws.Range("E1").Value = "cup_rental"

# Row 2 (customer 14045)
$ws.Range("D2").Value = 44866
$ws.Range("E2").ClearContents()
$ws.Range("G2").Value = 1
$ws.Range("H2").ClearContents()

# Row 3 (customer 14046)
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 0

# Row 4 (customer 14047)
$ws.Range("E4").Value = 989967
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 0

# Row 5 (customer 14048)
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 0

# Row 6 (customer 14049)
$ws.Range("E6").Value = 989969
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 1

# Move the selection like the author left it
$ws.Range("E5").Select()
